# Auto-generated Excel COM-interop script applying the DM-test correction
# (see commit: "Correcion a Diebold Mariano y revision de Cap1")
$wb = $excel.ActiveWorkbook

# --- Matriz_Resultados ---
$ws = $wb.Worksheets.Item("Matriz_Resultados")
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(7, 2).Value = 0
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 6).Value = 0

# --- P_valores ---
$ws = $wb.Worksheets.Item("P_valores")
$ws.Cells.Item(2, 3).Value = 0.002266742813555434
$ws.Cells.Item(2, 4).Value = 0.06187345410222744
$ws.Cells.Item(2, 5).Value = 0.32711061263428531
$ws.Cells.Item(2, 6).Value = 0.11478576708571619
$ws.Cells.Item(2, 7).Value = 0.01829374061026412
$ws.Cells.Item(2, 8).Value = 0.01782959888301572
$ws.Cells.Item(2, 9).Value = 0.003660262585102769
$ws.Cells.Item(2, 10).Value = 0.02112086554426473
$ws.Cells.Item(3, 2).Value = 0.002266742813555434
$ws.Cells.Item(3, 4).Value = 0.0000008614478133051051
$ws.Cells.Item(3, 5).Value = 0.0009254868924368775
$ws.Cells.Item(3, 6).Value = 0.001718523574625097
$ws.Cells.Item(3, 7).Value = 0.00005892973404475477
$ws.Cells.Item(3, 8).Value = 0.000005073847617165228
$ws.Cells.Item(3, 9).Value = 0.001143795705813133
$ws.Cells.Item(3, 10).Value = 0.000004804527816260062
$ws.Cells.Item(4, 2).Value = 0.06187345410222744
$ws.Cells.Item(4, 3).Value = 0.0000008614478133051051
$ws.Cells.Item(4, 5).Value = 0.01885580995447
$ws.Cells.Item(4, 6).Value = 0.04033026222480007
$ws.Cells.Item(4, 7).Value = 0.2633336196851277
$ws.Cells.Item(4, 8).Value = 0.12906253952169641
$ws.Cells.Item(4, 9).Value = 0.000005195122425583065
$ws.Cells.Item(4, 10).Value = 0.0002651291722188809
$ws.Cells.Item(5, 2).Value = 0.32711061263428531
$ws.Cells.Item(5, 3).Value = 0.0009254868924368775
$ws.Cells.Item(5, 4).Value = 0.01885580995447
$ws.Cells.Item(5, 6).Value = 0.50607790741186331
$ws.Cells.Item(5, 7).Value = 0.004759013923832267
$ws.Cells.Item(5, 8).Value = 0.004838830347924539
$ws.Cells.Item(5, 9).Value = 0.001553237744277292
$ws.Cells.Item(5, 10).Value = 0.37088247806760188
$ws.Cells.Item(6, 2).Value = 0.11478576708571619
$ws.Cells.Item(6, 3).Value = 0.001718523574625097
$ws.Cells.Item(6, 4).Value = 0.04033026222480007
$ws.Cells.Item(6, 5).Value = 0.50607790741186331
$ws.Cells.Item(6, 7).Value = 0.01139456376784831
$ws.Cells.Item(6, 8).Value = 0.01137617027982918
$ws.Cells.Item(6, 9).Value = 0.00274245977634191
$ws.Cells.Item(6, 10).Value = 0.06389298004274968
$ws.Cells.Item(7, 2).Value = 0.01829374061026412
$ws.Cells.Item(7, 3).Value = 0.00005892973404475477
$ws.Cells.Item(7, 4).Value = 0.2633336196851277
$ws.Cells.Item(7, 5).Value = 0.004759013923832267
$ws.Cells.Item(7, 6).Value = 0.01139456376784831
$ws.Cells.Item(7, 8).Value = 0.58535934186697114
$ws.Cells.Item(7, 9).Value = 0.0002261603920146182
$ws.Cells.Item(7, 10).Value = 0.00004052252392439115
$ws.Cells.Item(8, 2).Value = 0.01782959888301572
$ws.Cells.Item(8, 3).Value = 0.000005073847617165228
$ws.Cells.Item(8, 4).Value = 0.12906253952169641
$ws.Cells.Item(8, 5).Value = 0.004838830347924539
$ws.Cells.Item(8, 6).Value = 0.01137617027982918
$ws.Cells.Item(8, 7).Value = 0.58535934186697114
$ws.Cells.Item(8, 9).Value = 0.00004204109319827509
$ws.Cells.Item(8, 10).Value = 0.00002701894990431164
$ws.Cells.Item(9, 2).Value = 0.003660262585102769
$ws.Cells.Item(9, 3).Value = 0.001143795705813133
$ws.Cells.Item(9, 4).Value = 0.000005195122425583065
$ws.Cells.Item(9, 5).Value = 0.001553237744277292
$ws.Cells.Item(9, 6).Value = 0.00274245977634191
$ws.Cells.Item(9, 7).Value = 0.0002261603920146182
$ws.Cells.Item(9, 8).Value = 0.00004204109319827509
$ws.Cells.Item(9, 10).Value = 0.000010282254782501
$ws.Cells.Item(10, 2).Value = 0.02112086554426473
$ws.Cells.Item(10, 3).Value = 0.000004804527816260062
$ws.Cells.Item(10, 4).Value = 0.0002651291722188809
$ws.Cells.Item(10, 5).Value = 0.37088247806760188
$ws.Cells.Item(10, 6).Value = 0.06389298004274968
$ws.Cells.Item(10, 7).Value = 0.00004052252392439115
$ws.Cells.Item(10, 8).Value = 0.00002701894990431164
$ws.Cells.Item(10, 9).Value = 0.000010282254782501

# --- Estadisticos_DM ---
$ws = $wb.Worksheets.Item("Estadisticos_DM")
$ws.Cells.Item(2, 3).Value = 3.72406413379407608
$ws.Cells.Item(2, 4).Value = 2.0293997915935682
$ws.Cells.Item(2, 5).Value = -1.0154704353057129
$ws.Cells.Item(2, 6).Value = -1.6817163248867979
$ws.Cells.Item(2, 7).Value = 2.67017245057361086
$ws.Cells.Item(2, 8).Value = 2.68331313894356693
$ws.Cells.Item(2, 9).Value = 3.48250958281883882
$ws.Cells.Item(2, 10).Value = -2.59649594471495293
$ws.Cells.Item(3, 2).Value = -3.72406413379407608
$ws.Cells.Item(3, 4).Value = -8.32286617370647619
$ws.Cells.Item(3, 5).Value = -4.18021163988725775
$ws.Cells.Item(3, 6).Value = -3.86427026416571895
$ws.Cells.Item(3, 7).Value = -5.65939769739677256
$ws.Cells.Item(3, 8).Value = -7.13226377146948209
$ws.Cells.Item(3, 9).Value = -4.07163982205744457
$ws.Cells.Item(3, 10).Value = -7.16716433528292995
$ws.Cells.Item(4, 2).Value = -2.0293997915935682
$ws.Cells.Item(4, 3).Value = 8.32286617370647619
$ws.Cells.Item(4, 5).Value = -2.65468533077302604
$ws.Cells.Item(4, 6).Value = -2.25942627750447
$ws.Cells.Item(4, 7).Value = 1.1653895192674919
$ws.Cells.Item(4, 8).Value = 1.61294800313421804
$ws.Cells.Item(4, 9).Value = 7.11718061180051631
$ws.Cells.Item(4, 10).Value = -4.83392726635537606
$ws.Cells.Item(5, 2).Value = 1.0154704353057129
$ws.Cells.Item(5, 3).Value = 4.18021163988725775
$ws.Cells.Item(5, 4).Value = 2.65468533077302604
$ws.Cells.Item(5, 6).Value = 0.68246994283997442
$ws.Cells.Item(5, 7).Value = 3.35052266699115098
$ws.Cells.Item(5, 8).Value = 3.3421631220658119
$ws.Cells.Item(5, 9).Value = 3.91563172486782118
$ws.Cells.Item(5, 10).Value = -0.92450961959601008
$ws.Cells.Item(6, 2).Value = 1.6817163248867979
$ws.Cells.Item(6, 3).Value = 3.86427026416571895
$ws.Cells.Item(6, 4).Value = 2.25942627750447
$ws.Cells.Item(6, 5).Value = -0.68246994283997442
$ws.Cells.Item(6, 7).Value = 2.91087156528618918
$ws.Cells.Item(6, 8).Value = 2.91168878759655003
$ws.Cells.Item(6, 9).Value = 3.62789400550260988
$ws.Cells.Item(6, 10).Value = -2.01183694054372486
$ws.Cells.Item(7, 2).Value = -2.67017245057361086
$ws.Cells.Item(7, 3).Value = 5.65939769739677256
$ws.Cells.Item(7, 4).Value = -1.1653895192674919
$ws.Cells.Item(7, 5).Value = -3.35052266699115098
$ws.Cells.Item(7, 6).Value = -2.91087156528618918
$ws.Cells.Item(7, 8).Value = 0.55844893338993873
$ws.Cells.Item(7, 9).Value = 4.91894174832375608
$ws.Cells.Item(7, 10).Value = -5.87318538748175278
$ws.Cells.Item(8, 2).Value = -2.68331313894356693
$ws.Cells.Item(8, 3).Value = 7.13226377146948209
$ws.Cells.Item(8, 4).Value = -1.61294800313421804
$ws.Cells.Item(8, 5).Value = -3.3421631220658119
$ws.Cells.Item(8, 6).Value = -2.91168878759655003
$ws.Cells.Item(8, 7).Value = -0.55844893338993873
$ws.Cells.Item(8, 9).Value = 5.85202205602467718
$ws.Cells.Item(8, 10).Value = -6.10875868832475177
$ws.Cells.Item(9, 2).Value = -3.48250958281883882
$ws.Cells.Item(9, 3).Value = 4.07163982205744457
$ws.Cells.Item(9, 4).Value = -7.11718061180051631
$ws.Cells.Item(9, 5).Value = -3.91563172486782118
$ws.Cells.Item(9, 6).Value = -3.62789400550260988
$ws.Cells.Item(9, 7).Value = -4.91894174832375608
$ws.Cells.Item(9, 8).Value = -5.85202205602467718
$ws.Cells.Item(9, 10).Value = -6.68930433078311903
$ws.Cells.Item(10, 2).Value = 2.59649594471495293
$ws.Cells.Item(10, 3).Value = 7.16716433528292995
$ws.Cells.Item(10, 4).Value = 4.83392726635537606
$ws.Cells.Item(10, 5).Value = 0.92450961959601008
$ws.Cells.Item(10, 6).Value = 2.01183694054372486
$ws.Cells.Item(10, 7).Value = 5.87318538748175278
$ws.Cells.Item(10, 8).Value = 6.10875868832475177
$ws.Cells.Item(10, 9).Value = 6.68930433078311903

# --- Resumen ---
$ws = $wb.Worksheets.Item("Resumen")
$ws.Cells.Item(2, 2).Value = 6
$ws.Cells.Item(2, 4).Value = 2
$ws.Cells.Item(2, 5).Value = 75
$ws.Cells.Item(3, 2).Value = 4
$ws.Cells.Item(3, 4).Value = 3
$ws.Cells.Item(3, 5).Value = 50
$ws.Cells.Item(4, 1).Value = 'AV-MCPS'
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 4).Value = 5
$ws.Cells.Item(4, 5).Value = 12.5
$ws.Cells.Item(4, 6).Value = 0.69689255709081366
$ws.Cells.Item(5, 1).Value = 'MCPS'
$ws.Cells.Item(5, 2).Value = 1
$ws.Cells.Item(5, 4).Value = 5
$ws.Cells.Item(5, 5).Value = 12.5
$ws.Cells.Item(5, 6).Value = 0.70425960754298467
$ws.Cells.Item(6, 2).Value = 1
$ws.Cells.Item(6, 4).Value = 5
$ws.Cells.Item(6, 5).Value = 12.5
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 8
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(8, 4).Value = 7
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 8
